# Auto-generated Excel COM-interop script applying the Maduin_Profits.xlsx diff
# (workbook sheet names differ from the diff path; rows map 1:1 across sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1400.6
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5730
# Row 73
$ws.Range("H73").Value = 1400.6
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -5064
# Row 80
$ws.Range("H80").Value = 2804.5625
$ws.Range("J80").Value = 4000.375
$ws.Range("L80").Value = 12001.125
$ws.Range("N80").Value = -13997.125
# Row 83
$ws.Range("H83").Value = 2804.5625
$ws.Range("J83").Value = 4000.375
$ws.Range("L83").Value = 36003.375
$ws.Range("N83").Value = -45987.375
# Row 111
$ws.Range("H111").Value = 2796.5
$ws.Range("I111").Value = 2494.75
$ws.Range("K111").Value = 7484.25
$ws.Range("M111").Value = -4417.25
# Row 118
$ws.Range("H118").Value = 7500
$ws.Range("J118").Value = 7500
$ws.Range("L118").Value = 22500
$ws.Range("N118").Value = -25814
# Row 129
$ws.Range("H129").Value = 7565.6665
$ws.Range("I129").Value = 10598.5
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 31795.5
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = -26795.5
$ws.Range("N129").Value = -14500
# Row 138
$ws.Range("H138").Value = 7799
$ws.Range("J138").Value = 9530.888999999999
$ws.Range("L138").Value = 28592.667
$ws.Range("N138").Value = -38872.667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4291.42
$ws.Range("I32").Value = 3220.25
$ws.Range("K32").Value = 3220.25
$ws.Range("M32").Value = -2933.25
# Row 74
$ws.Range("H74").Value = 1558.2
$ws.Range("J74").Value = 1999.5
$ws.Range("L74").Value = 1999.5
$ws.Range("N74").Value = -3747.5
# Row 77
$ws.Range("H77").Value = 1558.2
$ws.Range("J77").Value = 1999.5
$ws.Range("L77").Value = 9997.5
$ws.Range("N77").Value = -18733.5
# Row 98
$ws.Range("H98").Value = 29498
$ws.Range("J98").Value = 29498
$ws.Range("L98").Value = 29498
$ws.Range("N98").Value = -35488
# Row 132
$ws.Range("H132").Value = 6109
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 7152.6
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 21457.8
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -26517.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 584.5
$ws.Range("I10").Value = 446.33334
$ws.Range("J10").Value = 999
$ws.Range("K10").Value = 446.33334
$ws.Range("L10").Value = 999
$ws.Range("M10").Value = -306.33334
$ws.Range("N10").Value = -1279
# Row 24
$ws.Range("H24").Value = 841
$ws.Range("I24").Value = 841
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 841
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -606
$ws.Range("N24").ClearContents()
# Row 86
$ws.Range("H86").Value = 2480.182
$ws.Range("I86").Value = 2075
$ws.Range("J86").Value = 3189.25
$ws.Range("K86").Value = 2075
$ws.Range("L86").Value = 3189.25
$ws.Range("M86").Value = -952
$ws.Range("N86").Value = -5435.25
# Row 89
$ws.Range("H89").Value = 2480.182
$ws.Range("I89").Value = 2075
$ws.Range("J89").Value = 3189.25
$ws.Range("K89").Value = 10375
$ws.Range("L89").Value = 15946.25
$ws.Range("M89").Value = -4759
$ws.Range("N89").Value = -27178.25
# Row 105
$ws.Range("H105").Value = 4284.5
$ws.Range("I105").Value = 1751
$ws.Range("J105").Value = 5129
$ws.Range("K105").Value = 1751
$ws.Range("L105").Value = 5129
$ws.Range("M105").Value = -4
$ws.Range("N105").Value = -8623
# Row 134
$ws.Range("H134").Value = 2150.2307
$ws.Range("I134").Value = 1873.4546
$ws.Range("K134").Value = 5620.3638
$ws.Range("M134").Value = -3085.3638

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 1292.5264
$ws.Range("I7").Value = 1072.6
$ws.Range("J7").Value = 1536.8889
$ws.Range("K7").Value = 1072.6
$ws.Range("L7").Value = 1536.8889
$ws.Range("M7").Value = -959.5999999999999
$ws.Range("N7").Value = -1762.8889
# Row 105
$ws.Range("H105").Value = 3067.5715
$ws.Range("I105").Value = 963.8570999999999
$ws.Range("K105").Value = 963.8570999999999
$ws.Range("M105").Value = 783.1429000000001
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 438.75
$ws.Range("I26").Value = 84.166664
$ws.Range("J26").Value = 1502.5
$ws.Range("K26").Value = 252.499992
$ws.Range("L26").Value = 4507.5
$ws.Range("M26").Value = 35.50000800000001
$ws.Range("N26").Value = -5083.5
# Row 132
$ws.Range("H132").Value = 2923.6428
$ws.Range("I132").Value = 1483.5
$ws.Range("K132").Value = 13351.5
$ws.Range("M132").Value = -10821.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6681386.5
$ws.Range("J70").Value = 12624.875
$ws.Range("L70").Value = 12624.875
$ws.Range("N70").Value = -13164.875
# Row 73
$ws.Range("H73").Value = 6681386.5
$ws.Range("J73").Value = 12624.875
$ws.Range("L73").Value = 12624.875
$ws.Range("N73").Value = -14496.875
# Row 80
$ws.Range("H80").Value = 2566.0625
$ws.Range("I80").Value = 2467.111
$ws.Range("J80").Value = 2693.2856
$ws.Range("K80").Value = 2467.111
$ws.Range("L80").Value = 2693.2856
$ws.Range("M80").Value = -1469.111
$ws.Range("N80").Value = -4689.2856
# Row 83
$ws.Range("H83").Value = 2566.0625
$ws.Range("I83").Value = 2467.111
$ws.Range("J83").Value = 2693.2856
$ws.Range("K83").Value = 12335.555
$ws.Range("L83").Value = 13466.428
$ws.Range("M83").Value = -7343.555
$ws.Range("N83").Value = -23450.428
# Row 113
$ws.Range("H113").Value = 1707.1428
$ws.Range("I113").Value = 1241.6666
$ws.Range("K113").Value = 1241.6666
$ws.Range("M113").Value = 928.3334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 888.86664
$ws.Range("I16").Value = 888.86664
$ws.Range("K16").Value = 888.86664
$ws.Range("M16").Value = -718.86664
# Row 20
$ws.Range("H20").Value = 6433
# Row 43
$ws.Range("H43").Value = 9999.666999999999
$ws.Range("I43").Value = 7999
$ws.Range("J43").Value = 11000
$ws.Range("K43").Value = 7999
$ws.Range("L43").Value = 11000
$ws.Range("M43").Value = -7806
$ws.Range("N43").Value = -11386
# Row 93
$ws.Range("H93").Value = 938.3
$ws.Range("I93").Value = 875.8889
$ws.Range("K93").Value = 875.8889
$ws.Range("M93").Value = 372.1111

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

